$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 858
$ws.Range("I12").Value = 810.3333
$ws.Range("K12").Value = 810.3333
$ws.Range("M12").Value = -640.3333
$ws.Range("H40").Value = 8330.583000000001
$ws.Range("I40").Value = 5999.4
$ws.Range("J40").Value = 9995.714
$ws.Range("K40").Value = 5999.4
$ws.Range("L40").Value = 9995.714
$ws.Range("M40").Value = -5824.4
$ws.Range("N40").Value = -10345.714
$ws.Range("H138").Value = 2548.4211
$ws.Range("J138").Value = 5076.625
$ws.Range("L138").Value = 15229.875
$ws.Range("N138").Value = -25509.875

$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 3954.3333
$ws.Range("I74").Value = 3348.8
$ws.Range("K74").Value = 3348.8
$ws.Range("M74").Value = -2474.8
$ws.Range("H77").Value = 3954.3333
$ws.Range("I77").Value = 3348.8
$ws.Range("K77").Value = 16744
$ws.Range("M77").Value = -12376
$ws.Range("H122").Value = 1666.3334
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H132").Value = 870.2692
$ws.Range("I132").Value = 859.5417
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 2578.6251
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -48.6251000000002
$ws.Range("N132").Value = -8057

$ws = $wb.Worksheets.Item(3)
$ws.Range("H64").Value = 773.5
$ws.Range("I64").Value = 718.6667
$ws.Range("J64").Value = 828.3333
$ws.Range("K64").Value = 718.6667
$ws.Range("L64").Value = 828.3333
$ws.Range("M64").Value = -493.6667
$ws.Range("N64").Value = -1278.3333
$ws.Range("H67").Value = 773.5
$ws.Range("I67").Value = 718.6667
$ws.Range("J67").Value = 828.3333
$ws.Range("K67").Value = 718.6667
$ws.Range("L67").Value = 828.3333
$ws.Range("M67").Value = 61.33330000000001
$ws.Range("N67").Value = -2388.3333

$ws = $wb.Worksheets.Item(4)
$ws.Range("H5").Value = 601.8
$ws.Range("I5").Value = 86.5
$ws.Range("J5").Value = 1374.75
$ws.Range("K5").Value = 86.5
$ws.Range("L5").Value = 1374.75
$ws.Range("M5").Value = 25.5
$ws.Range("N5").Value = -1598.75
$ws.Range("H6").Value = 404.5
$ws.Range("I6").Value = 404.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 404.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -291.5
$ws.Range("N6").ClearContents()
$ws.Range("H132").Value = 2241.6
$ws.Range("I132").Value = 2043.7894
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 6131.3682
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -3601.3682
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item(5)
$ws.Range("H9").Value = 1650
$ws.Range("I9").Value = 301
$ws.Range("K9").Value = 903
$ws.Range("M9").Value = -679
$ws.Range("H39").Value = 8606.637000000001
$ws.Range("J39").Value = 9327.299999999999
$ws.Range("L39").Value = 27981.9
$ws.Range("N39").Value = -28569.9
$ws.Range("H50").Value = 850
$ws.Range("I50").Value = 875
$ws.Range("J50").Value = 800
$ws.Range("K50").Value = 2625
$ws.Range("L50").Value = 2400
$ws.Range("M50").Value = -2144
$ws.Range("N50").Value = -3362
$ws.Range("H53").Value = 850
$ws.Range("I53").Value = 875
$ws.Range("J53").Value = 800
$ws.Range("K53").Value = 2625
$ws.Range("L53").Value = 2400
$ws.Range("M53").Value = -2144
$ws.Range("N53").Value = -3362
$ws.Range("H75").Value = 5781.8
$ws.Range("I75").Value = 333
$ws.Range("K75").Value = 999
$ws.Range("M75").Value = -1
$ws.Range("H78").Value = 5781.8
$ws.Range("I78").Value = 333
$ws.Range("K78").Value = 2997
$ws.Range("M78").Value = 1995
$ws.Range("H113").Value = 1243.091
$ws.Range("I113").Value = 815
$ws.Range("J113").Value = 1338.2222
$ws.Range("K113").Value = 2445
$ws.Range("L113").Value = 4014.6666
$ws.Range("M113").Value = -275
$ws.Range("N113").Value = -8354.6666
$ws.Range("H124").Value = 1111
$ws.Range("J124").Value = 1111
$ws.Range("L124").Value = 3333
$ws.Range("N124").Value = -13153
$ws.Range("H126").Value = 1800
$ws.Range("I126").Value = 1800
$ws.Range("K126").Value = 5400
$ws.Range("M126").Value = -460
$ws.Range("H132").Value = 3643.5
$ws.Range("I132").Value = 2924.7144
$ws.Range("K132").Value = 26322.4296
$ws.Range("M132").Value = -23792.4296

$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 6199.4
$ws.Range("I113").Value = 2684.4285
$ws.Range("J113").Value = 9275
$ws.Range("K113").Value = 2684.4285
$ws.Range("L113").Value = 9275
$ws.Range("M113").Value = -514.4285
$ws.Range("N113").Value = -13615
$ws.Range("H122").Value = 2124.5715
$ws.Range("I122").Value = 1557.6
$ws.Range("K122").Value = 4672.799999999999
$ws.Range("M122").Value = -2222.799999999999
$ws.Range("H133").Value = 135000
$ws.Range("J133").Value = 135000
$ws.Range("L133").Value = 135000
$ws.Range("N133").Value = -145120

$ws = $wb.Worksheets.Item(7)
$ws.Range("H55").Value = 1298.1875
$ws.Range("I55").Value = 1704.8
$ws.Range("K55").Value = 1704.8
$ws.Range("M55").Value = -1531.8
$ws.Range("H122").Value = 2272.9443
$ws.Range("I122").Value = 2365.6924
$ws.Range("K122").Value = 7097.0772
$ws.Range("M122").Value = -4647.0772
$ws.Range("H123").Value = 79999
$ws.Range("J123").Value = 79999
$ws.Range("L123").Value = 79999
$ws.Range("N123").Value = -89799
$ws.Range("H132").Value = 3499.5715
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 3750
$ws.Range("I136").Value = 3750
$ws.Range("K136").Value = 11250
$ws.Range("M136").Value = -8700

$ws = $wb.Worksheets.Item(8)
$ws.Range("H3").Value = 4500
$ws.Range("I3").Value = 4500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4386
$ws.Range("N3").ClearContents()
$ws.Range("H11").Value = 1999
$ws.Range("J11").Value = 1999
$ws.Range("L11").Value = 1999
$ws.Range("N11").Value = -2283
$ws.Range("H29").Value = 4950
$ws.Range("J29").Value = 4950
$ws.Range("L29").Value = 4950
$ws.Range("N29").Value = -5530
$ws.Range("H32").Value = 14000.5
$ws.Range("I32").Value = 8001
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 8001
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -7684
$ws.Range("N32").Value = -20634
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H122").Value = 2144.1738
$ws.Range("I122").Value = 1980.55
$ws.Range("K122").Value = 5941.65
$ws.Range("M122").Value = -3491.65
